$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 374.53333
$ws.Range("I19").Value2 = 266.14285
$ws.Range("K19").Value2 = 266.14285
$ws.Range("M19").Value2 = -91.14285000000001

$ws.Range("H40").Value2 = 4349.3887
$ws.Range("J40").Value2 = 4487.5884
$ws.Range("L40").Value2 = 4487.5884
$ws.Range("N40").Value2 = -4837.5884

$ws.Range("H53").Value2 = 2344.6667
$ws.Range("I53").Value2 = 39.333332
$ws.Range("K53").Value2 = 39.333332
$ws.Range("M53").Value2 = 597.666668

$ws.Range("H92").Value2 = 832.28125
$ws.Range("I92").Value2 = 613.9643
$ws.Range("K92").Value2 = 613.9643
$ws.Range("M92").Value2 = 634.0357

$ws.Range("H97").Value2 = 3937.25
$ws.Range("I97").Value2 = 2249.5
$ws.Range("K97").Value2 = 6748.5
$ws.Range("M97").Value2 = -6252.5

$ws.Range("H116").Value2 = 2700.5
$ws.Range("I116").Value2 = 2561.5
$ws.Range("K116").Value2 = 2561.5
$ws.Range("M116").Value2 = 880.5

$ws.Range("H125").Value2 = 2107.625
$ws.Range("I125").Value2 = 893.5
$ws.Range("J125").Value2 = 5750
$ws.Range("K125").Value2 = 8041.5
$ws.Range("L125").Value2 = 51750
$ws.Range("M125").Value2 = -5581.5
$ws.Range("N125").Value2 = -56670

$ws.Range("H127").Value2 = 814.38464
$ws.Range("J127").Value2 = 1914
$ws.Range("L127").Value2 = 5742
$ws.Range("N127").Value2 = -15662

$ws.Range("H131").Value2 = 3546.52
$ws.Range("I131").Value2 = 3652.6667
$ws.Range("J131").Value2 = 999
$ws.Range("K131").Value2 = 10958.0001
$ws.Range("L131").Value2 = 2997
$ws.Range("M131").Value2 = -5918.000100000001
$ws.Range("N131").Value2 = -13077

$ws.Range("H132").Value2 = 15875001
$ws.Range("I132").Value2 = 17243400
$ws.Range("J132").Value2 = 1559.4
$ws.Range("K132").Value2 = 51730200
$ws.Range("L132").Value2 = 4678.200000000001
$ws.Range("M132").Value2 = -51727670
$ws.Range("N132").Value2 = -9738.200000000001

$ws.Range("H133").Value2 = 41999.668
$ws.Range("J133").Value2 = 41999.668
$ws.Range("L133").Value2 = 41999.668
$ws.Range("N133").Value2 = -52119.668

$ws.Range("H136").Value2 = 44000
$ws.Range("J136").Value2 = 44000
$ws.Range("L136").Value2 = 44000
$ws.Range("N136").Value2 = -54200

$ws.Range("H138").Value2 = 264872.97
$ws.Range("I138").Value2 = 807.1739
$ws.Range("J138").Value2 = 669773.9
$ws.Range("K138").Value2 = 2421.5217
$ws.Range("L138").Value2 = 2009321.7
$ws.Range("M138").Value2 = 2718.4783
$ws.Range("N138").Value2 = -2019601.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 7236.857
$ws.Range("J61").Value2 = 7583.212
$ws.Range("L61").Value2 = 7583.212
$ws.Range("N61").Value2 = -8007.212

$ws.Range("H74").Value2 = 2691.024
$ws.Range("I74").Value2 = 820.7692
$ws.Range("J74").Value2 = 5730.1875
$ws.Range("K74").Value2 = 820.7692
$ws.Range("L74").Value2 = 5730.1875
$ws.Range("M74").Value2 = 53.23080000000004
$ws.Range("N74").Value2 = -7478.1875

$ws.Range("H77").Value2 = 2691.024
$ws.Range("I77").Value2 = 820.7692
$ws.Range("J77").Value2 = 5730.1875
$ws.Range("K77").Value2 = 4103.846
$ws.Range("L77").Value2 = 28650.9375
$ws.Range("M77").Value2 = 264.1540000000005
$ws.Range("N77").Value2 = -37386.9375

$ws.Range("H132").Value2 = 1932.4263
$ws.Range("I132").Value2 = 1895.1637
$ws.Range("K132").Value2 = 5685.4911
$ws.Range("M132").Value2 = -3155.4911

$ws.Range("H136").Value2 = 7236.857
$ws.Range("J136").Value2 = 7583.212
$ws.Range("L136").Value2 = 22749.636
$ws.Range("N136").Value2 = -27849.636

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 2026831
$ws.Range("I86").Value2 = 4453638
$ws.Range("J86").Value2 = 4491.8335
$ws.Range("K86").Value2 = 4453638
$ws.Range("L86").Value2 = 4491.8335
$ws.Range("M86").Value2 = -4452515
$ws.Range("N86").Value2 = -6737.8335

$ws.Range("H89").Value2 = 2026831
$ws.Range("I89").Value2 = 4453638
$ws.Range("J89").Value2 = 4491.8335
$ws.Range("K89").Value2 = 22268190
$ws.Range("L89").Value2 = 22459.1675
$ws.Range("M89").Value2 = -22262574
$ws.Range("N89").Value2 = -33691.1675

$ws.Range("H105").Value2 = 6010.1577
$ws.Range("I105").Value2 = 5137.909
$ws.Range("K105").Value2 = 5137.909
$ws.Range("M105").Value2 = -3390.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value2 = 34505.25
$ws.Range("I103").Value2 = 34505.25
$ws.Range("K103").Value2 = 34505.25
$ws.Range("M103").Value2 = -33333.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 527.19354
$ws.Range("I12").Value2 = 586.8333
$ws.Range("K12").Value2 = 1760.4999
$ws.Range("M12").Value2 = -1587.4999

$ws.Range("H33").Value2 = 587.1667
$ws.Range("J33").Value2 = 697.2
$ws.Range("L33").Value2 = 4183.200000000001
$ws.Range("N33").Value2 = -4749.200000000001

$ws.Range("H40").Value2 = 1352.75
$ws.Range("I40").Value2 = 1540.3334
$ws.Range("J40").Value2 = 790
$ws.Range("K40").Value2 = 6161.3336
$ws.Range("L40").Value2 = 3160
$ws.Range("M40").Value2 = -6092.3336
$ws.Range("N40").Value2 = -3298

$ws.Range("H132").Value2 = 966.53845
$ws.Range("J132").Value2 = 1300
$ws.Range("L132").Value2 = 11700
$ws.Range("N132").Value2 = -16760

$ws.Range("H136").Value2 = 6718.75
$ws.Range("J136").Value2 = 6718.75
$ws.Range("L136").Value2 = 20156.25
$ws.Range("N136").Value2 = -30356.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value2 = 48111.25
$ws.Range("J46").Value2 = 49555.715
$ws.Range("L46").Value2 = 49555.715
$ws.Range("N46").Value2 = -49867.715

$ws.Range("H80").Value2 = 4739.727
$ws.Range("I80").Value2 = 4762.4287
$ws.Range("K80").Value2 = 4762.4287
$ws.Range("M80").Value2 = -3764.4287

$ws.Range("H83").Value2 = 4739.727
$ws.Range("I83").Value2 = 4762.4287
$ws.Range("K83").Value2 = 23812.1435
$ws.Range("M83").Value2 = -18820.1435

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 5121.737
$ws.Range("J40").Value2 = 4685.4287
$ws.Range("L40").Value2 = 4685.4287
$ws.Range("N40").Value2 = -4957.4287

$ws.Range("H108").Value2 = 100000
$ws.Range("J108").Value2 = 100000
$ws.Range("L108").Value2 = 100000
$ws.Range("N108").Value2 = -107680

$ws.Range("H122").Value2 = 3809.5417
$ws.Range("I122").Value2 = 3466.4546
$ws.Range("J122").Value2 = 4099.846
$ws.Range("K122").Value2 = 10399.3638
$ws.Range("L122").Value2 = 12299.538
$ws.Range("M122").Value2 = -7949.363799999999
$ws.Range("N122").Value2 = -17199.538

$ws.Range("H132").Value2 = 3308.7
$ws.Range("I132").Value2 = 3425.0454
$ws.Range("J132").Value2 = 3166.5
$ws.Range("K132").Value2 = 10275.1362
$ws.Range("L132").Value2 = 9499.5
$ws.Range("M132").Value2 = -7745.136200000001
$ws.Range("N132").Value2 = -14559.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 63588.707
$ws.Range("J81").Value2 = 9662.5
$ws.Range("L81").Value2 = 19325
$ws.Range("N81").Value2 = -21447

$ws.Range("H84").Value2 = 63588.707
$ws.Range("J84").Value2 = 9662.5
$ws.Range("L84").Value2 = 96625
$ws.Range("N84").Value2 = -107233

$ws.Range("H100").Value2 = 1140.5
$ws.Range("I100").Value2 = 1220.1428
$ws.Range("J100").Value2 = 861.75
$ws.Range("K100").Value2 = 2440.2856
$ws.Range("L100").Value2 = 1723.5
$ws.Range("M100").Value2 = -1899.2856
$ws.Range("N100").Value2 = -2805.5

$ws.Range("H113").Value2 = 3971982.5
$ws.Range("I113").Value2 = 4389928
$ws.Range("K113").Value2 = 13169784
$ws.Range("M113").Value2 = -13167614

$ws.Range("H120").Value2 = 86666.664
$ws.Range("J120").Value2 = 86666.664
$ws.Range("L120").Value2 = 86666.664
$ws.Range("N120").Value2 = -96342.664

$ws.Range("H124").Value2 = 100171.4
$ws.Range("J124").Value2 = 100171.4
$ws.Range("L124").Value2 = 100171.4
$ws.Range("N124").Value2 = -109991.4

$ws.Range("H126").Value2 = 2539.7
$ws.Range("I126").Value2 = 2223.1765
$ws.Range("J126").Value2 = 4333.3335
$ws.Range("K126").Value2 = 6669.529500000001
$ws.Range("L126").Value2 = 13000.0005
$ws.Range("M126").Value2 = -4199.529500000001
$ws.Range("N126").Value2 = -17940.0005

$ws.Range("H132").Value2 = 4030.4
$ws.Range("I132").Value2 = 4413.125
$ws.Range("K132").Value2 = 13239.375
$ws.Range("M132").Value2 = -10709.375

$ws.Range("H136").Value2 = 7058.5293
$ws.Range("I136").Value2 = 7651.2905
$ws.Range("J136").Value2 = 933.3333
$ws.Range("K136").Value2 = 22953.8715
$ws.Range("L136").Value2 = 2799.9999
$ws.Range("M136").Value2 = -20403.8715
$ws.Range("N136").Value2 = -7899.9999
